$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# -- Fill in the newly finalized statistics for "Health Sciences" (row 2)
#    and "Physical Sciences" (row 3), columns C:F (Papers/Authors x 2+/11+ etc.)
$ws.Range("C2").Value = 8758846
$ws.Range("D2").Value = 17195880
$ws.Range("E2").Value = 3441064
$ws.Range("F2").Value = 8705266

$ws.Range("C3").Value = 12096908
$ws.Range("D3").Value = 15366570
$ws.Range("E3").Value = 5182626
$ws.Range("F3").Value = 11557780

# -- Update the view state left behind by the author: scrolled one column
#    to the right (so column B is first visible column), zoomed to 170%,
#    with F3 as the active selection.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 170
$ws.Range("F3").Select()
